$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round the row-5 sensor readings to 2 decimal places ("custom accuracy") ---
$ws.Range("B5").Value = 13.32
$ws.Range("C5").Value = 9.94
$ws.Range("D5").Value = 0.97
$ws.Range("E5").Value = 29.27
$ws.Range("F5").Value = 23.68
$ws.Range("G5").Value = 10.16
$ws.Range("H5").Value = 38.27
$ws.Range("I5").Value = 16.22
$ws.Range("J5").Value = 7.34
$ws.Range("K5").Value = 10.4
$ws.Range("L5").Value = 11.73
$ws.Range("M5").Value = 12.54
$ws.Range("N5").Value = 3.54
$ws.Range("O5").Value = 10.53
$ws.Range("P5").Value = 14.81
$ws.Range("Q5").Value = 9.02
$ws.Range("R5").Value = 0.44
$ws.Range("S5").Value = 0.55
$ws.Range("T5").Value = 152.88
$ws.Range("U5").Value = 29.33
$ws.Range("V5").Value = 9.72
$ws.Range("W5").Value = 19.56
$ws.Range("X5").Value = 10.35
$ws.Range("Y5").Value = 1.37
$ws.Range("Z5").Value = 19.43
$ws.Range("AA5").Value = 8.58
$ws.Range("AB5").Value = 7.66
$ws.Range("AC5").Value = 8.98
$ws.Range("AD5").Value = 12.37
$ws.Range("AE5").Value = 0.48
$ws.Range("AF5").Value = 34.65
$ws.Range("AG5").Value = 5.4
$ws.Range("AH5").Value = 12.15

# --- Drop the now-redundant last data row (row 6) -> dimension shrinks to A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- Narrow a handful of data columns from 8 to 7 characters wide ---
# (6.15 lands in the same "snap to 7 characters" bucket Excel's pixel-rounding
# uses for an 8->7 ColumnWidth edit on this Calibri 11 sheet.)
$ws.Columns.Item(7).ColumnWidth = 6.15
$ws.Columns.Item(11).ColumnWidth = 6.15
$ws.Columns.Item(15).ColumnWidth = 6.15
$ws.Columns.Item(24).ColumnWidth = 6.15
